# Personal Journal #5 | Client/Professor Meeting Journal Update
#
# Appends a new journal entry ("#5 - Backend setup and progress") after the
# existing "#4" entry, following the same pattern used by the earlier
# entries in the document: a blank separator paragraph, a heading
# paragraph, and a body paragraph that begins with a tab character.

$d = $word.ActiveDocument

$enDash = [char]0x2013

$heading = "#5 " + $enDash + " Backend setup and progress"
$body = "After creating basic navigation with react, I decided to begin work on the backend. I thought this would be a good idea before diving deeper into reacts hook, state, and redux. So far it has been slightly easier than I expected. I was able to build the frontend quicker than anticipated, and I guessed that the backend setup would be much more difficult in comparison. However, I feel I am making good progress. After finishing up the API endpoints in the backend, my next step is to tie it all together back in the frontend. "

# Build the WordprocessingML for the three new paragraphs:
#   1. an empty separator paragraph (matches the blank "<w:p/>" that sits
#      between every entry in this journal)
#   2. the "#5 - Backend setup and progress" heading paragraph
#   3. the body paragraph, whose single run starts with a real <w:tab/>
#      element (not a literal tab character) followed by the journal text
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$bodyXml = "<w:body>" +
    "<w:p/>" +
    "<w:p><w:r><w:t>" + $heading + "</w:t></w:r></w:p>" +
    "<w:p><w:r><w:tab/><w:t xml:space=`"preserve`">" + $body + "</w:t></w:r></w:p>" +
    "</w:body>"

$packageXml = '<?xml version="1.0" standalone="yes"?>' +
    '<?mso-application progid="Word.Document"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="' + $wNs + '">' + $bodyXml + '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

# Move to the very end of the document and insert the new content there.
$r = $d.Content
$r.Collapse(0)
$r.InsertXML($packageXml) | Out-Null
